$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" row (row 26) and, after the shift, the "SC 92" row
# (which becomes row 27) -- everything below shifts up accordingly.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Apply the remaining value changes (new missing-data pattern) using the
# post-deletion row numbers.
$ws.Range("F4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("D6").Value = -14.2
$ws.Range("F7").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E11").Value = -7.9
$ws.Range("F17").Value = 17.78
$ws.Range("D19").Value = -15.5
$ws.Range("E19").ClearContents()
$ws.Range("D21").ClearContents()
$ws.Range("D23").Value = -13.9
$ws.Range("E23").Value = -7
$ws.Range("F24").Value = 16.78
$ws.Range("E25").Value = -7.1
$ws.Range("C26").ClearContents()
$ws.Range("C27").Value = 10
$ws.Range("D27").ClearContents()
$ws.Range("E27").ClearContents()
$ws.Range("F27").Value = 17
$ws.Range("F28").Value = 17.44
$ws.Range("C29").ClearContents()
$ws.Range("D29").Value = -13
$ws.Range("E29").ClearContents()
$ws.Range("F29").ClearContents()
$ws.Range("E30").Value = -5.7
$ws.Range("F30").ClearContents()
$ws.Range("F32").ClearContents()
$ws.Range("E33").Value = -10.7
